$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Set-ParagraphXml($paraIndex, $innerXml) {
    $p = $d.Paragraphs.Item($paraIndex).Range
    $xml = '<?xml version="1.0"?><w:document ' + $wNs + '><w:body>' + $innerXml + '</w:body></w:document>'
    $p.InsertXML($xml)
}

# 1) "Depuis les travaux de la philosophe Annah Arendt..." -> "Hannah", proofErr removed,
#    text split into "Ha" + "nnah Arendt..." runs.
Set-ParagraphXml 3 '<w:p><w:r><w:t xml:space="preserve">Depuis les travaux de la philosophe </w:t></w:r><w:r><w:t>Ha</w:t></w:r><w:r><w:t>nnah Arendt (1951) et du politiste Carl Joachim Friedrich (1953), il désigne également le nazisme et le stalinisme.</w:t></w:r></w:p>'

# 2) "Adolf Hitler, le Führer, ... Joseph Staline Vodj" -> remove proofErr around "Vodj" (stays italic).
Set-ParagraphXml 14 '<w:p><w:r><w:t xml:space="preserve">Adolf Hitler, le </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>Führer</w:t></w:r><w:r><w:t xml:space="preserve">, Benito Mussolini </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>Duce</w:t></w:r><w:r><w:t xml:space="preserve">, Joseph Staline </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>Vodj</w:t></w:r></w:p>'

# 3) "Führer, Duce, Vodj = le chef" -> merge three runs into one plain run (no italics), proofErr removed.
Set-ParagraphXml 15 '<w:p><w:r><w:t>Führer, Duce, Vodj = le chef</w:t></w:r></w:p>'

# 4) "1937 Les japonais commettent le sac de Nanka." -> merge into one run, proofErr removed.
Set-ParagraphXml 22 '<w:p><w:r><w:t>1937 Les japonais commettent le sac de Nanka.</w:t></w:r></w:p>'

# 5) "Les chemises noires (Camicie nere)" -> merge italic runs into one, proofErr removed.
Set-ParagraphXml 24 '<w:p><w:r><w:t>Les chemises noires (</w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>Camicie nere</w:t></w:r><w:r><w:t>)</w:t></w:r></w:p>'

# 6) "Ce discours a été tenu dans la brasserie (Hofbraühaus) à Munich." -> proofErr removed only.
Set-ParagraphXml 32 '<w:p><w:r><w:t>Ce discours a été tenu dans la brasserie (</w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>Hofbraühaus</w:t></w:r><w:r><w:t>) à Munich.</w:t></w:r></w:p>'

# 7) "Affiche Ein Volk, Ein reich, Ein Führer" -> merge runs, proofErr removed.
Set-ParagraphXml 36 '<w:p><w:r><w:t>Affiche Ein Volk, Ein reich, Ein Führer</w:t></w:r></w:p>'

# 8) "Le Volk est figuré par la famille allemande idéale..." -> merge runs, proofErr removed.
Set-ParagraphXml 38 '<w:p><w:r><w:t>Le Volk est figuré par la famille allemande idéale (blonde), composée d’un couple de trois enfants.</w:t></w:r></w:p>'

# 9) "Union soviétique, lénino-marxisme" -> merge runs, proofErr removed.
Set-ParagraphXml 55 '<w:p><w:r><w:t>Union soviétique, lénino-marxisme</w:t></w:r></w:p>'

# 10) "...Hitler m'a dit, H. Rauschning, 1979)" -> merge trailing runs, proofErr removed, keep italic book title run.
Set-ParagraphXml 58 '<w:p><w:r><w:t xml:space="preserve">C’est avec la jeunesse […] dompter la peur » (issue du livre, </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>Hitler m’a dit</w:t></w:r><w:r><w:t>, H. Rauschning, 1979)</w:t></w:r></w:p>'

# 11) Bold the two "Regarder le film" paragraphs.
Set-ParagraphXml 60 '<w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Regarder le film : Le garçon au pyjama rayé.</w:t></w:r></w:p>'
Set-ParagraphXml 61 '<w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Regarder le film : Jojo Rabbit</w:t></w:r></w:p>'

Write-Output "done"
